# Updated: st 13. 01. 2021
# Refresh the Slovakia COVID daily-stats sheet:
#  - correct AgTests (H) / AgPosit (I) figures for several existing rows
#  - append the new day's row (314)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to previously-reported AgTests (H) / AgPosit (I) values ---
$ws.Range("H266").Value = 12549
$ws.Range("I266").Value = 662
$ws.Range("H267").Value = 13130
$ws.Range("I267").Value = 779
$ws.Range("H268").Value = 14021
$ws.Range("I268").Value = 687
$ws.Range("H287").Value = 57796
$ws.Range("I287").Value = 3928
$ws.Range("H288").Value = 56583
$ws.Range("I288").Value = 3990
$ws.Range("H289").Value = 64863
$ws.Range("I289").Value = 3715
$ws.Range("H292").Value = 81928
$ws.Range("I292").Value = 7258
$ws.Range("H293").Value = 82860
$ws.Range("I293").Value = 5859
$ws.Range("H294").Value = 92000
$ws.Range("I294").Value = 5099
$ws.Range("H299").Value = 65000
$ws.Range("I299").Value = 6821
$ws.Range("H300").Value = 70834
$ws.Range("I300").Value = 6930
$ws.Range("H301").Value = 69789
$ws.Range("I301").Value = 5528
$ws.Range("H302").Value = 72746
$ws.Range("I302").Value = 5290
$ws.Range("H306").Value = 70319
$ws.Range("I306").Value = 7147
$ws.Range("H307").Value = 72394
$ws.Range("I307").Value = 6233
$ws.Range("H309").Value = 56571
$ws.Range("I309").Value = 3910
$ws.Range("H310").Value = 89093
$ws.Range("I310").Value = 5333
$ws.Range("H311").Value = 32044
$ws.Range("I311").Value = 1249
$ws.Range("H312").Value = 37894
$ws.Range("I312").Value = 1130
$ws.Range("H313").Value = 71472
$ws.Range("I313").Value = 3450

# --- Append the new daily row (314) ---
$ws.Range("A314").Value = 44208
$ws.Range("A314").NumberFormat = "yyyy-mm-dd"
$ws.Range("B314").Value = 215055
$ws.Range("C314").Value = 157028
$ws.Range("D314").Value = 54864
$ws.Range("E314").Value = 15440
$ws.Range("F314").Value = 3576
$ws.Range("G314").Value = 3163
$ws.Range("H314").Value = 62078
$ws.Range("I314").Value = 3188
